$d = $word.ActiveDocument

# --- Paragraph 1: the "**ID__AFFARS_..._ID**" placeholder paragraph ---
$p1 = $d.Paragraphs.Item(1)

# Add a paragraph border (w:pBdr) with 5pt spacing on all four sides -
# mirrors the w:space="5" pBdr already present on the third paragraph.
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5

# w:ind w:left goes from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Update the placeholder id text in the first run.
$d.Content.Find.Execute("**ID__AFFARS_5350_topic_5__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_5350_102_1_70__ID**", 2)

# Remove the trailing single-space run left over at the end of paragraph 1.
$p1 = $d.Paragraphs.Item(1)
$tail = $p1.Range.Duplicate
$tail.SetRange($tail.End - 2, $tail.End - 1)
if ($tail.Text -eq " ") {
    $tail.Delete()
}
